# Fruta / hortaliza, semanal
# Insert a new weekly record as row 230 (pushing existing rows 230-294
# down to 231-295) on the single data sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 230:294 down one row, creating a blank row 230.
$ws.Rows(230).Insert()

# Populate the new row 230 with the new weekly observation.
$ws.Cells.Item(230, 1).Value2  = 6
$ws.Cells.Item(230, 2).Value2  = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(230, 3).Value2  = "Metropolitana"
$ws.Cells.Item(230, 4).Value2  = 44588
$ws.Cells.Item(230, 5).Value2  = 13
$ws.Cells.Item(230, 6).Value2  = "Fruta"
$ws.Cells.Item(230, 7).Value2  = 100101
$ws.Cells.Item(230, 8).Value2  = "Berries"
$ws.Cells.Item(230, 9).Value2  = 100101001
$ws.Cells.Item(230, 10).Value2 = "Arándano (blue)"
$ws.Cells.Item(230, 11).Value2 = "Sin especificar"
$ws.Cells.Item(230, 12).Value2 = "Primera"
$ws.Cells.Item(230, 13).Value2 = 3500
$ws.Cells.Item(230, 14).Value2 = 4000
$ws.Cells.Item(230, 15).Value2 = 4000
$ws.Cells.Item(230, 16).Value2 = 4000
$ws.Cells.Item(230, 17).Value2 = "`$/bandeja 2 kilos"
$ws.Cells.Item(230, 18).Value2 = "Provincia de Linares"
$ws.Cells.Item(230, 19).Value2 = 2000
$ws.Cells.Item(230, 20).Value2 = 2

# Match the date-formatted style used by the other rows in column D.
$ws.Cells.Item(230, 4).NumberFormat = $ws.Cells.Item(231, 4).NumberFormat
